# Generate Report for Handback
# Adds a new handback record (85292b46-8ed0-4aa3-815a-34da459a2008.md)
# as row 4 to the "Overview", "zh-cn" and "de-de" tables/worksheets.

$wb = $excel.ActiveWorkbook

$fileGuid = "85292b46-8ed0-4aa3-815a-34da459a2008"
$fileName = "$fileGuid.md"
$pathName = "e2e\$fileGuid.md"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOverview = $rowOverview.Range

$rngOverview.Cells.Item(1, 1).Value = $fileName
$rngOverview.Cells.Item(1, 2).Value = $pathName
$rngOverview.Cells.Item(1, 3).Value = ".md"
$rngOverview.Cells.Item(1, 5).Value = "Handed back: in sync with en-US"
$rngOverview.Cells.Item(1, 6).Value = "Handed back: in sync with en-US"
$rngOverview.Cells.Item(1, 7).Value = "2016-08-21 22:53:56"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/751115f0bdfee1cb3c0da18bd921f396c49b3025/e2e/$fileName", $null, $null, $pathName) | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range

$zhCnXlf = "$fileGuid.751115f0bdfee1cb3c0da18bd921f396c49b3025.zh-cn.xlf"

$rngZhCn.Cells.Item(1, 1).Value = $fileName
$rngZhCn.Cells.Item(1, 2).Value = ".md"
$rngZhCn.Cells.Item(1, 3).Value = "Handed back: in sync with en-US"
$rngZhCn.Cells.Item(1, 4).Value = "e2e"
$rngZhCn.Cells.Item(1, 5).Value = "ht"
$rngZhCn.Cells.Item(1, 6).Value = "True"
$rngZhCn.Cells.Item(1, 7).Value = $zhCnXlf
$rngZhCn.Cells.Item(1, 8).Value = "2016-08-21 22:53:52"
$rngZhCn.Cells.Item(1, 9).Value = $fileName
$rngZhCn.Cells.Item(1, 10).Value = $zhCnXlf
$rngZhCn.Cells.Item(1, 11).Value = "2016-08-21 22:54:14"
$rngZhCn.Cells.Item(1, 12).Value = ""
$rngZhCn.Cells.Item(1, 13).Value = "True"
$rngZhCn.Cells.Item(1, 14).Value = ""
$rngZhCn.Cells.Item(1, 15).Value = "False"
$rngZhCn.Cells.Item(1, 16).Value = ""

$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/751115f0bdfee1cb3c0da18bd921f396c49b3025/e2e/$fileName", $null, $null, $fileName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/751115f0bdfee1cb3c0da18bd921f396c49b3025/e2e/$fileName", $null, $null, $fileName) | Out-Null

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range

$deDeXlf = "$fileGuid.751115f0bdfee1cb3c0da18bd921f396c49b3025.de-de.xlf"

$rngDeDe.Cells.Item(1, 1).Value = $fileName
$rngDeDe.Cells.Item(1, 2).Value = ".md"
$rngDeDe.Cells.Item(1, 3).Value = "Handed back: in sync with en-US"
$rngDeDe.Cells.Item(1, 4).Value = "e2e"
$rngDeDe.Cells.Item(1, 5).Value = "ht"
$rngDeDe.Cells.Item(1, 6).Value = "True"
$rngDeDe.Cells.Item(1, 7).Value = $deDeXlf
$rngDeDe.Cells.Item(1, 8).Value = "2016-08-21 22:53:56"
$rngDeDe.Cells.Item(1, 9).Value = $fileName
$rngDeDe.Cells.Item(1, 10).Value = $deDeXlf
$rngDeDe.Cells.Item(1, 11).Value = "2016-08-21 22:54:20"
$rngDeDe.Cells.Item(1, 12).Value = ""
$rngDeDe.Cells.Item(1, 13).Value = "True"
$rngDeDe.Cells.Item(1, 14).Value = ""
$rngDeDe.Cells.Item(1, 15).Value = "False"
$rngDeDe.Cells.Item(1, 16).Value = ""

$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/751115f0bdfee1cb3c0da18bd921f396c49b3025/e2e/$fileName", $null, $null, $fileName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/751115f0bdfee1cb3c0da18bd921f396c49b3025/e2e/$fileName", $null, $null, $fileName) | Out-Null
